$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (E1)
$ws.Range("C2").Value = 22
$ws.Range("D2").Value = 497
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 0.04
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.08
$ws.Range("I2").Value = 0.03
$ws.Range("J2").Value = 0.5
$ws.Range("K2").Value = 0.06

# Row 3 (E2)
$ws.Range("F3").Value = 0.03
$ws.Range("H3").Value = 0.05
$ws.Range("I3").Value = 0.02
$ws.Range("K3").Value = 0.03

# Row 4 (E3)
$ws.Range("F4").Value = 0.03
$ws.Range("G4").Value = 0.71
$ws.Range("H4").Value = 0.06
$ws.Range("I4").Value = 0.02
$ws.Range("J4").Value = 0.73
$ws.Range("K4").Value = 0.04

# Row 5 (E4)
$ws.Range("F5").Value = 0.04
$ws.Range("G5").Value = 0.86
$ws.Range("H5").Value = 0.07
$ws.Range("I5").Value = 0.02
$ws.Range("J5").Value = 0.83
$ws.Range("K5").Value = 0.04

# Row 6 (E5)
$ws.Range("F6").Value = 0.07
$ws.Range("G6").Value = 0.37
$ws.Range("H6").Value = 0.11
$ws.Range("I6").Value = 0.05
$ws.Range("J6").Value = 0.47
$ws.Range("K6").Value = 0.09

# Row 7 (E6)
$ws.Range("F7").Value = 0.06
$ws.Range("G7").Value = 0.69
$ws.Range("H7").Value = 0.11
$ws.Range("I7").Value = 0.04
$ws.Range("J7").Value = 0.67
$ws.Range("K7").Value = 0.07

# Row 8 (E7)
$ws.Range("F8").Value = 0.06
$ws.Range("G8").Value = 0.47
$ws.Range("H8").Value = 0.1
$ws.Range("I8").Value = 0.04
$ws.Range("K8").Value = 0.07
